$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Tracking" labels to disambiguate tracking dimension 1 vs 2
$ws.Range("C5").Value = "Tracking-1: {{track_name}}"
$ws.Range("C8").Value = "Tracking-2: {{track2_name}}"
$ws.Range("E10").Value = "Tracking-1"

# Update the selected cell shown when the sheet is opened
$ws.Range("C15").Select()
